$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D retains its original text formatting so numeric-looking
# price strings (e.g. "1.00", "7.09") are not reinterpreted as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# --- Price (D) / Volume(1h) (E) updates ---
$ws.Range("D2").Value = "69.298.57"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "3.663.63"
$ws.Range("E3").Value = "  -0.68%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "643.12"
$ws.Range("E5").Value = "  -5.68%  "
$ws.Range("D6").Value = "158.98"
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -1.09%  "
$ws.Range("D10").Value = "7.09"
$ws.Range("E10").Value = "  -0.62%  "
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("E12").Value = "  -0.33%  "
$ws.Range("D13").Value = "4.282.46"
$ws.Range("E13").Value = "  -0.66%  "
$ws.Range("D14").Value = "32.48"
$ws.Range("E14").Value = "  -0.06%  "
$ws.Range("D15").Value = "3.646.37"
$ws.Range("E15").Value = "  -0.89%  "
$ws.Range("D16").Value = "69.296.63"
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("E17").Value = "  +0.40%  "
$ws.Range("E18").Value = "  -0.50%  "
$ws.Range("D19").Value = "6.45"
$ws.Range("E19").Value = "  -0.49%  "
$ws.Range("D20").Value = "466.41"
$ws.Range("E20").Value = "  -0.87%  "
$ws.Range("D21").Value = "9.81"
$ws.Range("E21").Value = "  -1.13%  "
$ws.Range("D22").Value = "0.643"
$ws.Range("E22").Value = "  -1.79%  "
$ws.Range("D23").Value = "79.35"
$ws.Range("E23").Value = "  -0.84%  "
$ws.Range("D24").Value = "3.812.40"
$ws.Range("E24").Value = "  -0.61%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("E26").Value = "  +0.91%  "
$ws.Range("D27").Value = "10.82"
$ws.Range("E27").Value = "  -1.14%  "
$ws.Range("E28").Value = "  -2.45%  "
$ws.Range("D29").Value = "2.61"
$ws.Range("E29").Value = "  -3.11%  "
$ws.Range("D30").Value = "1.70"
$ws.Range("E30").Value = "  -2.75%  "
$ws.Range("D31").Value = "1.99"
$ws.Range("E31").Value = "  -0.26%  "
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("E33").Value = "  -1.20%  "
$ws.Range("D34").Value = "6.44"
$ws.Range("E34").Value = "  -2.88%  "
$ws.Range("E35").Value = "  +3.48%  "
$ws.Range("D36").Value = "3.654.45"
$ws.Range("E36").Value = "  -0.61%  "
$ws.Range("E37").Value = "  +1.53%  "
$ws.Range("D39").Value = "5.88"
$ws.Range("E39").Value = "  -6.27%  "
$ws.Range("D40").Value = "177.91"
$ws.Range("E40").Value = "  +4.61%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").Value = "0.0896"
$ws.Range("E42").Value = "  -1.31%  "
$ws.Range("D43").Value = "2.17"
$ws.Range("E43").Value = "  -3.54%  "
$ws.Range("D44").Value = "0.924"
$ws.Range("E44").Value = "  -2.08%  "
$ws.Range("D45").Value = "46.51"
$ws.Range("E45").Value = "  -2.40%  "
$ws.Range("E46").Value = "  +0.13%  "

# --- Rows 47-51 reordering (Coin / Link / Price / Volume(1h)) ---
$ws.Range("B47").Value = "FLOKI"
$ws.Range("C47").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D47").Value = "0.000269"
$ws.Range("E47").Value = "  -2.89%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").Value = "1.25"
$ws.Range("E48").Value = "  -3.70%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "26.83"
$ws.Range("E49").Value = "  -5.04%  "
$ws.Range("B50").Value = "SuiNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D50").Value = "1.06"
$ws.Range("E50").Value = "  -5.29%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").Value = "7.81"
$ws.Range("E51").Value = "  +0.18%  "
